# cryptos.xlsx refresh - GitHub Actions scheduled update
#
# Updates the Price (D) and Volume(1h) (E) columns for the crypto list with
# freshly-scraped figures, and re-orders four rows (35-38) whose relative
# ranking changed (LidoDAOToken / ImmutableX / BinanceUSD / WEMIXToken),
# rewriting their Coin/Link/Price/Volume cells in place.
#
# Price values are stored as plain text in this sheet (e.g. "36.495.96",
# "0.657", "1.00") rather than numbers, so list-like separators and trailing
# zeros survive untouched. A handful of the new prices look like an ordinary
# decimal (e.g. "47.64"); Excel would silently re-interpret a bare Value
# assignment like that as a Number, so for those specific cells the new text
# is entered through Formula with a leading apostrophe (the standard Excel
# "force text" input prefix) to keep the cell a text cell, matching the rest
# of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.495.96"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "2.101.61"
$ws.Range("E3").Value = "  +9.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("E5").Value = "  +1.03%  "

# Row 6
$ws.Range("D6").Formula = "'0.657"
$ws.Range("E6").Value = "  -6.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Formula = "'47.64"
$ws.Range("E8").Value = "  +8.09%  "

# Row 9
$ws.Range("D9").Formula = "'59.59"
$ws.Range("E9").Value = "  +2.85%  "

# Row 10
$ws.Range("D10").Formula = "'0.375"
$ws.Range("E10").Value = "  +1.63%  "

# Row 11
$ws.Range("D11").Formula = "'0.0745"
$ws.Range("E11").Value = "  -2.38%  "

# Row 12
$ws.Range("E12").Value = "  +0.04%  "

# Row 13
$ws.Range("D13").Value = "2.410.29"
$ws.Range("E13").Value = "  +9.97%  "

# Row 14
$ws.Range("D14").Formula = "'14.48"
$ws.Range("E14").Value = "  -0.28%  "

# Row 15
$ws.Range("D15").Formula = "'0.827"
$ws.Range("E15").Value = "  +2.32%  "

# Row 16
$ws.Range("D16").Value = "2.105.65"
$ws.Range("E16").Value = "  +9.65%  "

# Row 17
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").Value = "36.517.65"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").Formula = "'72.86"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0829"
$ws.Range("E20").Value = "  -3.62%  "

# Row 21
$ws.Range("D21").Formula = "'13.31"
$ws.Range("E21").Value = "  -0.30%  "

# Row 22
$ws.Range("D22").Formula = "'239.80"
$ws.Range("E22").Value = "  -4.06%  "

# Row 23
$ws.Range("D23").Formula = "'5.14"
$ws.Range("E23").Value = "  -1.20%  "

# Row 25
$ws.Range("D25").Formula = "'2.46"
$ws.Range("E25").Value = "  -6.60%  "

# Row 26
$ws.Range("D26").Formula = "'171.08"
$ws.Range("E26").Value = "  +1.85%  "

# Row 27
$ws.Range("D27").Formula = "'21.43"
$ws.Range("E27").Value = "  +14.58%  "

# Row 28
$ws.Range("D28").Formula = "'9.13"
$ws.Range("E28").Value = "  +3.47%  "

# Row 29
$ws.Range("E29").Value = "  -9.84%  "

# Row 30
$ws.Range("D30").Formula = "'28.47"
$ws.Range("E30").Value = "  +59.98%  "

# Row 31
$ws.Range("E31").Value = "  -4.72%  "

# Row 32
$ws.Range("D32").Formula = "'0.0615"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33
$ws.Range("D33").Formula = "'4.43"
$ws.Range("E33").Value = "  -3.31%  "

# Row 34
$ws.Range("D34").Formula = "'0.0911"
$ws.Range("E34").Value = "  +2.18%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Formula = "'0.955"
$ws.Range("E35").Value = "  +8.67%  "

# Row 36
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Formula = "'1.00"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Formula = "'1.88"
$ws.Range("E37").Value = "  -3.62%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Formula = "'2.32"
$ws.Range("E38").Value = "  +13.85%  "

# Row 39
$ws.Range("E39").Value = "  -6.13%  "

# Row 40
$ws.Range("D40").Formula = "'1.33"
$ws.Range("E40").Value = "  -12.48%  "

# Row 41
$ws.Range("D41").Formula = "'1.17"
$ws.Range("E41").Value = "  +6.33%  "

# Row 42
$ws.Range("D42").Formula = "'0.0222"
$ws.Range("E42").Value = "  -1.92%  "

# Row 43
$ws.Range("D43").Formula = "'97.58"
$ws.Range("E43").Value = "  -8.75%  "

# Row 44
$ws.Range("D44").Formula = "'2.74"
$ws.Range("E44").Value = "  -4.25%  "

# Row 45
$ws.Range("D45").Formula = "'16.06"
$ws.Range("E45").Value = "  -7.44%  "

# Row 46
$ws.Range("D46").Value = "1.341.76"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("D47").Formula = "'0.0843"
$ws.Range("E47").Value = "  +3.42%  "

# Row 48
$ws.Range("D48").Formula = "'7.05"
$ws.Range("E48").Value = "  +10.35%  "

# Row 49
$ws.Range("D49").Value = "2.297.73"
$ws.Range("E49").Value = "  +9.85%  "

# Row 50
$ws.Range("E50").Value = "  +1.35%  "

# Row 51
$ws.Range("E51").Value = "  -5.67%  "
